$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.006.18"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "2.644.54"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.68"
$ws.Range("E5").Value = "  +3.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.58"
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  +0.81%  "
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("E10").Value = "  +5.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.352"
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "3.104.44"
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("D14").Value = "60.993.10"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.91"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("E16").Value = "  +2.63%  "
$ws.Range("D17").Value = "2.650.43"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "354.28"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.68"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.70"
$ws.Range("E23").Value = "  +2.16%  "
$ws.Range("E24").Value = "  +2.48%  "
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.994"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").Value = "0.0₃0866"
$ws.Range("E27").Value = "  +3.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.40"
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +7.19%  "
$ws.Range("E31").Value = "  +4.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.50"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.52"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.16"
$ws.Range("E34").Value = "  +4.29%  "
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.925"
$ws.Range("E36").Value = "  +9.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.895"
$ws.Range("E37").Value = "  +1.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "308.69"
$ws.Range("E38").Value = "  +4.67%  "
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.643"
$ws.Range("E41").Value = "  +3.75%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.102"
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0564"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.84"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.02"
$ws.Range("E46").Value = "  +3.01%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0239"
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.34"
$ws.Range("E48").Value = "  +8.26%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.34"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.990.35"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.84"
$ws.Range("E51").Value = "  +2.41%  "
